# Add "Norway" and "Poland" market test-data sheets to the gallery workbook,
# mirroring the existing per-country sheets (e.g. "Turkey") in layout/styles.

$wb = $excel.ActiveWorkbook

# --- Norway -----------------------------------------------------------
$template = $wb.Worksheets.Item("Turkey")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$template.Copy($null, $lastSheet) | Out-Null

$norway = $wb.Worksheets.Item($wb.Worksheets.Count)
$norway.Name = "Norway"
# Write the ticket id before the market label so the shared-string table
# picks up "NGC-2931/T3063" ahead of "Norway Market" (matches source order).
$norway.Range("B4").Value = "NGC-2931/T3063"
$norway.Range("B2").Value = "Norway Market"
$norway.Range("E23").Select() | Out-Null

# --- Poland -------------------------------------------------------------
$template2 = $wb.Worksheets.Item("Turkey")
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$template2.Copy($null, $lastSheet2) | Out-Null

$poland = $wb.Worksheets.Item($wb.Worksheets.Count)
$poland.Name = "Poland"
$poland.Range("B4").Value = "NGC-2920/T3106"
$poland.Range("B2").Value = "Poland Market"
$poland.Range("E23").Select() | Out-Null

# The workbook's active tab ends up on "Norway" (not the most recently
# added "Poland"), matching the source edit.
$norway.Activate() | Out-Null
